# Updated symbol list on Fri Dec 23 11:28:09 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) and rank/trend label (column E) values
# for the crypto ticker sheet to match the latest scrape. A handful of
# coins also changed rank/position (rows 41-43 shuffled between
# BKEXToken / CEJI / KickToken), so their Coin (B), Link (C), Price (D)
# and Data (E) cells are rewritten together to keep each row consistent.
#
# Note: column D stores prices as text (e.g. "246.14"), not numbers, so
# a leading apostrophe is used to force Excel to keep these as text
# instead of auto-converting them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'246.14"

# Row 3 - OKB
$ws.Range("D3").Value = "'22.02"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.426"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.05837"

# Row 8 - MXToken
$ws.Range("D8").Value = "'0.8079"

# Row 9 - FTXToken
$ws.Range("D9").Value = "'0.9654"

# Row 10
$ws.Range("D10").Value = "'0.1429"

# Row 11
$ws.Range("D11").Value = "'0.07456"

# Row 12 - LiechtensteinCryptoassetsExchange (now "Best in 24h")
$ws.Range("D12").Value = "'0.03321"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCXBestin24h"

# Row 13
$ws.Range("D13").Value = "'0.03018"

# Row 14
$ws.Range("D14").Value = "'4.167"

# Row 15
$ws.Range("D15").Value = "'0.09401"

# Row 16
$ws.Range("D16").Value = "'0.001596"

# Row 17
$ws.Range("D17").Value = "'0.04821"

# Row 18 - One (now "Worst in 24h")
$ws.Range("D18").Value = "'0.0005892"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# Row 19
$ws.Range("D19").Value = "'0.006134"

# Row 20
$ws.Range("D20").Value = "'0.004107"

# Row 21
$ws.Range("D21").Value = "'0.0009978"

# Row 22
$ws.Range("D22").Value = "'0.0001500"

# Row 24
$ws.Range("D24").Value = "'2.217"

# Row 26
$ws.Range("D26").Value = "'0.1297"

# Row 40 - IDEX
$ws.Range("D40").Value = "'0.03876"

# Row 41 - now KickToken (was BKEXToken)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006680"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 - now BKEXToken (was CEJI)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1077"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 - now CEJI (was KickToken)
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002540"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.006686"

# Row 45 - CoinLion
$ws.Range("D45").Value = "'0.00005621"

# Row 47 - CoinbaseStockToken (no longer "Best in 24h")
$ws.Range("D47").Value = "'0.4202"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

# Row 48 - BOLO
$ws.Range("D48").Value = "'0.1464"

# Row 49 - CryptobidCoin
$ws.Range("D49").Value = "'0.00002101"
